$d = $word.ActiveDocument

# The document ends with:
#   "...purrr::map()."                                   <- keep
#   "References:"                                         <- remove (whole paragraph)
#   "johnramey. <hyperlink>...in R</hyperlink>."           <- remove (whole paragraph)
#   ""                                                     <- keep (trailing empty paragraph)
#
# Remove the "References:" paragraph and the following reference-list
# paragraph in their entirety (including their paragraph marks), so the
# "...purrr::map()." paragraph is immediately followed by the trailing
# empty paragraph.

$refPara = $null
$nextPara = $null

foreach ($p in $d.Paragraphs) {
    if ($refPara -eq $null) {
        if ($p.Range.Text.TrimEnd("`r", "`a") -eq "References:") {
            $refPara = $p
        }
    } else {
        $nextPara = $p
        break
    }
}

if ($refPara -ne $null -and $nextPara -ne $null) {
    $delRange = $d.Range($refPara.Range.Start, $nextPara.Range.End)
    $delRange.Delete()
}
